$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at C and D (existing C,D shift to E,F; existing E shifts to G)
$ws.Range("C1:D1").EntireColumn.Insert()

# Header row
$ws.Range("C1").Value = "modelo"
$ws.Range("D1").Value = "politica"

# Data rows: set modelo (C), politica (D), full (E), tipo (F), link (G)
$ws.Range("C2").Value = "FONTE 60A LITE"
$ws.Range("D2").Value = "Igual"
$ws.Range("E2").Value = "NA"
$ws.Range("F2").Value = "classico"
$ws.Range("G2").Value = "https://www.mercadolivre.com.br/jfa-fonte-carregador-storm-lite-60a-3000-w-preto/p/MLB23456525?pdp_filters=seller_id:1492722080#searchVariation=MLB23456525&position=3&search_layout=stack&type=product&tracking_id=e7177dfe-d068-49c3-91db-5e16a39e7498"

$ws.Range("C3").Value = "FONTE 200 BOB"
$ws.Range("D3").Value = "Igual"
$ws.Range("E3").Value = "NA"
$ws.Range("F3").Value = "classico"
$ws.Range("G3").Value = "https://www.mercadolivre.com.br/fonte-automotiva-jfa-storm-200a-bob-carregador-automatico-bivolt-cor-bob-200a-jfa/p/MLB24834408?pdp_filters=seller_id:1492722080#searchVariation=MLB24834408&position=4&search_layout=stack&type=product&tracking_id=e7177dfe-d068-49c3-91db-5e16a39e7498"

$ws.Range("C4").Value = "FONTE 70A"
$ws.Range("D4").Value = "Igual"
$ws.Range("E4").Value = "NA"
$ws.Range("F4").Value = "classico"
$ws.Range("G4").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-70a-bivolt-com-medidor-cca/p/MLB21455208?pdp_filters=seller_id:1492722080#searchVariation=MLB21455208&position=5&search_layout=stack&type=product&tracking_id=e7177dfe-d068-49c3-91db-5e16a39e7498"

$ws.Range("C5").Value = "FONTE 90 BOB"
$ws.Range("D5").Value = "Igual"
$ws.Range("E5").Value = "NA"
$ws.Range("F5").Value = "classico"
$ws.Range("G5").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-bob-storm-90a-bivolt-automatico-cor-preto/p/MLB21562641?pdp_filters=seller_id:1492722080#searchVariation=MLB21562641&position=6&search_layout=stack&type=product&tracking_id=e7177dfe-d068-49c3-91db-5e16a39e7498"

$ws.Range("C6").Value = "FONTE 120A LITE"
$ws.Range("D6").Value = "Igual"
$ws.Range("E6").Value = "NA"
$ws.Range("F6").Value = "classico"
$ws.Range("G6").Value = "https://www.mercadolivre.com.br/fonte-carregador-automotivo-jfa-120a-storm-lite-12v-bivolt-cor-preto/p/MLB23998473?pdp_filters=seller_id:1492722080#searchVariation=MLB23998473&position=7&search_layout=stack&type=product&tracking_id=e7177dfe-d068-49c3-91db-5e16a39e7498"

$ws.Range("C7").Value = "FONTE 60A"
$ws.Range("D7").Value = "Igual"
$ws.Range("E7").Value = "NA"
$ws.Range("F7").Value = "classico"
$ws.Range("G7").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-60a-bivolt-storm-com-medidor-cca/p/MLB21320712?pdp_filters=seller_id:1492722080#searchVariation=MLB21320712&position=8&search_layout=stack&type=product&tracking_id=e7177dfe-d068-49c3-91db-5e16a39e7498"

$ws.Range("C8").Value = "Sem Modelo"
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = "NA"
$ws.Range("F8").Value = "classico"
$ws.Range("G8").Value = "https://www.mercadolivre.com.br/fonte-carregador-automotivo-jfa-40a-lite-storm-slim-bivolt-cor-preto/p/MLB33435981?pdp_filters=seller_id:1492722080#searchVariation=MLB33435981&position=9&search_layout=stack&type=product&tracking_id=e7177dfe-d068-49c3-91db-5e16a39e7498"

$ws.Range("C9").Value = "FONTE 40A"
$ws.Range("D9").Value = "Igual"
$ws.Range("E9").Value = "NA"
$ws.Range("F9").Value = "classico"
$ws.Range("G9").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-storm-40a-bivolt-12v-cor-preto/p/MLB22569833?pdp_filters=seller_id:1492722080#searchVariation=MLB22569833&position=10&search_layout=stack&type=product&tracking_id=e7177dfe-d068-49c3-91db-5e16a39e7498"

$ws.Range("C10").Value = "FONTE 200 MONO"
$ws.Range("D10").Value = "Igual"
$ws.Range("E10").Value = "NA"
$ws.Range("F10").Value = "classico"
$ws.Range("G10").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-200a-storm-voltimetro-digital-mono-220v-cor-preto/p/MLB24006449?pdp_filters=seller_id:1492722080#searchVariation=MLB24006449&position=11&search_layout=stack&type=product&tracking_id=e7177dfe-d068-49c3-91db-5e16a39e7498"

$ws.Range("C11").Value = "FONTE 200A LITE"
$ws.Range("D11").Value = "Igual"
$ws.Range("E11").Value = "NA"
$ws.Range("F11").Value = "classico"
$ws.Range("G11").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-200a-lite-storm-slim-bivolt-cor-azul/p/MLB24154371?pdp_filters=seller_id:1492722080#searchVariation=MLB24154371&position=12&search_layout=stack&type=product&tracking_id=e7177dfe-d068-49c3-91db-5e16a39e7498"

$ws.Range("C12").Value = "Sem Modelo"
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "classico"
$ws.Range("G12").Value = "https://www.mercadolivre.com.br/fonte-carregador-automotivo-jfa-200a-storm-lite-mono-220v-cor-azul/p/MLB30464905?pdp_filters=seller_id:1492722080#searchVariation=MLB30464905&position=13&search_layout=stack&type=product&tracking_id=e7177dfe-d068-49c3-91db-5e16a39e7498"

$ws.Range("C13").Value = "FONTE 200A"
$ws.Range("D13").Value = "Igual"
$ws.Range("E13").Value = "NA"
$ws.Range("F13").Value = "premium"
$ws.Range("G13").Value = "https://www.mercadolivre.com.br/fonte-carregador-automotiva-storm-sci-redline-jfa-200a-slim-cor-preto/p/MLB26091477?pdp_filters=seller_id:1492722080#searchVariation=MLB26091477&position=1&search_layout=stack&type=product&tracking_id=e7177dfe-d068-49c3-91db-5e16a39e7498"

$ws.Range("C14").Value = "FONTE 200 BOB"
$ws.Range("D14").Value = "Igual"
$ws.Range("E14").Value = "NA"
$ws.Range("F14").Value = "premium"
$ws.Range("G14").Value = "https://produto.mercadolivre.com.br/MLB-3643553599-fonte-carregador-automotivo-jfa-bob-storm-200a-bivolt-_JM?searchVariation=180217356804#searchVariation%3D180217356804%26position%3D14%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3De7177dfe-d068-49c3-91db-5e16a39e7498"

$ws.Range("C15").Value = "FONTE 70A LITE"
$ws.Range("D15").Value = "Igual"
$ws.Range("E15").Value = "NA"
$ws.Range("F15").Value = "classico"
$ws.Range("G15").Value = "https://produto.mercadolivre.com.br/MLB-4589672834-fonte-carregador-automotivo-jfa-70a-storm-lite-12v-bivolt-_JM#position%3D15%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3De7177dfe-d068-49c3-91db-5e16a39e7498"

$ws.Range("C16").Value = "Sem Modelo"
$ws.Range("D16").Value = ""
$ws.Range("E16").Value = "NA"
$ws.Range("F16").Value = "premium"
$ws.Range("G16").Value = "https://produto.mercadolivre.com.br/MLB-3548356679-fonte-carregador-automotivo-jfa-200a-storm-lite-mono-220v-_JM#position%3D16%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3De7177dfe-d068-49c3-91db-5e16a39e7498"

$ws.Range("C17").Value = "FONTE 200A LITE"
$ws.Range("D17").Value = "Igual"
$ws.Range("E17").Value = "NA"
$ws.Range("F17").Value = "premium"
$ws.Range("G17").Value = "https://produto.mercadolivre.com.br/MLB-3651101141-fonte-carregador-automotivo-jfa-200a-storm-lite-12v-bivolt-_JM#position%3D17%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3De7177dfe-d068-49c3-91db-5e16a39e7498"

$ws.Range("C18").Value = "FONTE 70A"
$ws.Range("D18").Value = "Igual"
$ws.Range("E18").Value = "NA"
$ws.Range("F18").Value = "premium"
$ws.Range("G18").Value = "https://produto.mercadolivre.com.br/MLB-3643619483-fonte-automotiva-jfa-storm-70a-bivolt-com-medidor-cca-e-sci-_JM#position%3D18%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3De7177dfe-d068-49c3-91db-5e16a39e7498"

$ws.Range("C19").Value = "FONTE 60A"
$ws.Range("D19").Value = "Igual"
$ws.Range("E19").Value = "NA"
$ws.Range("F19").Value = "premium"
$ws.Range("G19").Value = "https://produto.mercadolivre.com.br/MLB-3643752475-fonte-automotiva-jfa-storm-60a-bivolt-com-medidor-cca-e-sci-_JM#position%3D19%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3De7177dfe-d068-49c3-91db-5e16a39e7498"

$ws.Range("C20").Value = "FONTE 60A LITE"
$ws.Range("D20").Value = "Igual"
$ws.Range("E20").Value = "NA"
$ws.Range("F20").Value = "premium"
$ws.Range("G20").Value = "https://produto.mercadolivre.com.br/MLB-3646752187-fonte-carregador-automotivo-jfa-60a-storm-lite-12v-bivolt-_JM#position%3D20%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3De7177dfe-d068-49c3-91db-5e16a39e7498"

$ws.Range("C21").Value = "FONTE 200 MONO"
$ws.Range("D21").Value = "Igual"
$ws.Range("E21").Value = "NA"
$ws.Range("F21").Value = "premium"
$ws.Range("G21").Value = "https://produto.mercadolivre.com.br/MLB-4589396544-fonte-carregador-jfa-200a-storm-voltimetro-digital-mono-220v-_JM#position%3D21%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3De7177dfe-d068-49c3-91db-5e16a39e7498"

$ws.Range("C22").Value = "FONTE 120A LITE"
$ws.Range("D22").Value = "Igual"
$ws.Range("E22").Value = "NA"
$ws.Range("F22").Value = "premium"
$ws.Range("G22").Value = "https://produto.mercadolivre.com.br/MLB-3641319051-fonte-carregador-automotivo-jfa-120a-storm-lite-12v-bivolt-_JM#position%3D22%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3De7177dfe-d068-49c3-91db-5e16a39e7498"
